$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title for the new Prolog results table (merged F12:I12), mirrors F1's style
$ws.Range("F12:I12").Merge() | Out-Null
$ws.Range("F12").Value = "Prolog, primes until 10000, upper limit: 2^1000000"
$ws.Range("F12:I12").HorizontalAlignment = -4108

# Header row (13) mirrors row 13's A:D headers
$ws.Range("F13").Value = "Processes"
$ws.Range("G13").Value = "Execution time (s)"
$ws.Range("H13").Value = "Speedup"
$ws.Range("I13").Value = "Efficiency"

# Data rows 14-21
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 283.81

$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 178.55
$ws.Range("H15").Formula = "=G14/G15"
$ws.Range("I15").Formula = "=H15/F15"

$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 153
$ws.Range("H16").Formula = "=G14/G16"
$ws.Range("I16").Formula = "=H16/F16"

$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 151.28
$ws.Range("H17").Formula = "=G14/G17"
$ws.Range("I17").Formula = "=H17/F17"

$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 150.56
$ws.Range("H18").Formula = "=G14/G18"
$ws.Range("I18").Formula = "=H18/F18"

$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 151.19
$ws.Range("H19").Formula = "=G14/G19"
$ws.Range("I19").Formula = "=H19/F19"

$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 150.63
$ws.Range("H20").Formula = "=G14/G20"
$ws.Range("I20").Formula = "=H20/F20"

$ws.Range("F21").Value = 8
$ws.Range("G21").Value = 151.03
$ws.Range("H21").Formula = "=G14/G21"
$ws.Range("I21").Formula = "=H21/F21"

# Update the selection to match target state
$ws.Range("H18").Select() | Out-Null
